$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.009.26'
$ws.Range('E2').Value = '  +4.75%  '
$ws.Range('D3').Value = '2.622.63'
$ws.Range('E3').Value = '  +5.56%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '605.78'
$ws.Range('E5').Value = '  +3.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '180.81'
$ws.Range('E6').Value = '  +3.81%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +2.05%  '
$ws.Range('D9').Value = '2.621.62'
$ws.Range('E9').Value = '  +5.58%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  +14.88%  '
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('E12').Value = '  +3.52%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.03'
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('E14').Value = '  +5.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.65'
$ws.Range('E15').Value = '  +5.50%  '
$ws.Range('E16').Value = '  +7.49%  '
$ws.Range('D17').Value = '70.587.39'
$ws.Range('E17').Value = '  +4.30%  '
$ws.Range('D18').Value = '2.619.27'
$ws.Range('E18').Value = '  +5.60%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '380.40'
$ws.Range('E19').Value = '  +9.92%  '
$ws.Range('E20').Value = '  +7.31%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.51'
$ws.Range('E21').Value = '  +6.66%  '
$ws.Range('E22').Value = '  +1.51%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.96'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.44'
$ws.Range('E24').Value = '  +6.52%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.84'
$ws.Range('E26').Value = '  +9.86%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.64'
$ws.Range('E27').Value = '  +9.27%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '0.0₃0951'
$ws.Range('E30').Value = '  +7.19%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '525.80'
$ws.Range('E31').Value = '  +5.62%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.02'
$ws.Range('E32').Value = '  +4.05%  '
$ws.Range('E33').Value = '  +6.72%  '
$ws.Range('E34').Value = '  +4.24%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '164.07'
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.120'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  +4.72%  '
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.87'
$ws.Range('E40').Value = '  +8.46%  '
$ws.Range('E41').Value = '  +5.36%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.03'
$ws.Range('E43').Value = '  +5.76%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.60'
$ws.Range('E44').Value = '  +9.63%  '
$ws.Range('E45').Value = '  +2.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '40.12'
$ws.Range('E46').Value = '  +3.92%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '153.70'
$ws.Range('E47').Value = '  +4.40%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0273'
$ws.Range('E48').Value = '  +8.35%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.65'
$ws.Range('E49').Value = '  +3.85%  '
$ws.Range('E50').Value = '  +4.41%  '
